# Generate Report for Handoff
# Updates the localization-status report: the zh-cn / de-de rows move from
# "In Translation" to "Ready for handoff", the handoff timestamps are
# refreshed, and the Status columns are widened to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -----------------------
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refresh the handoff-generation timestamps ------------------------------
$wsZhCn.Range("H2").Value     = "2016-08-23 17:03:34"   # zh-cn Latest Handoff Datetime
$wsOverview.Range("G2").Value = "2016-08-23 17:03:38"   # Latest HO Xliff Generate Date
$wsDeDe.Range("H2").Value     = "2016-08-23 17:03:38"   # de-de Latest Handoff Datetime

# --- Widen the Status columns to fit "Ready for handoff" --------------------
# ColumnWidth is quantized by Excel to whole pixels, so 16.35 is the input
# that lands on the pixel step closest to the target character width.
$wsOverview.Columns.Item(5).ColumnWidth = 16.35   # zh-cn column
$wsOverview.Columns.Item(6).ColumnWidth = 16.35   # de-de column
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.35   # Status column
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.35   # Status column
